$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link swap (row 34 <-> row 35) ---
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'

# --- Price (D) and Volume(1h) (E) updates ---
# Use a text-number-format round trip so numeric-looking strings
# (e.g. "1.002") are stored as literal text, matching the source data,
# then restore the base style so no stray number-format style sticks.
$updates = @(
    @{ Ref = 'D2'; Value = '26.883.11' },
    @{ Ref = 'E2'; Value = '  -0.37%  ' },
    @{ Ref = 'D3'; Value = '1.859.34' },
    @{ Ref = 'E3'; Value = '  +0.03%  ' },
    @{ Ref = 'D4'; Value = '1.002' },
    @{ Ref = 'E4'; Value = '  -0.34%  ' },
    @{ Ref = 'D5'; Value = '304.49' },
    @{ Ref = 'E5'; Value = '  -0.50%  ' },
    @{ Ref = 'D6'; Value = '1.001' },
    @{ Ref = 'E6'; Value = '  -0.07%  ' },
    @{ Ref = 'D7'; Value = '0.5044' },
    @{ Ref = 'E7'; Value = '  -1.09%  ' },
    @{ Ref = 'D8'; Value = '0.3642' },
    @{ Ref = 'D9'; Value = '0.07161' },
    @{ Ref = 'E9'; Value = '  +0.87%  ' },
    @{ Ref = 'D10'; Value = '0.8912' },
    @{ Ref = 'E10'; Value = '  +0.80%  ' },
    @{ Ref = 'D11'; Value = '20.67' },
    @{ Ref = 'E11'; Value = '  +0.84%  ' },
    @{ Ref = 'D12'; Value = '1.878.57' },
    @{ Ref = 'E12'; Value = '  +0.93%  ' },
    @{ Ref = 'D13'; Value = '0.07476' },
    @{ Ref = 'E13'; Value = '  -0.97%  ' },
    @{ Ref = 'D14'; Value = '93.73' },
    @{ Ref = 'E14'; Value = '  +5.95%  ' },
    @{ Ref = 'D15'; Value = '5.226' },
    @{ Ref = 'E15'; Value = '  -1.46%  ' },
    @{ Ref = 'D16'; Value = '1.002' },
    @{ Ref = 'E16'; Value = '  -0.40%  ' },
    @{ Ref = 'D17'; Value = '0.000008491' },
    @{ Ref = 'E17'; Value = '  +0.55%  ' },
    @{ Ref = 'D18'; Value = '14.19' },
    @{ Ref = 'E18'; Value = '  +0.86%  ' },
    @{ Ref = 'D19'; Value = '1.001' },
    @{ Ref = 'E19'; Value = '  -0.22%  ' },
    @{ Ref = 'D20'; Value = '26.948.67' },
    @{ Ref = 'E20'; Value = '  -0.41%  ' },
    @{ Ref = 'D21'; Value = '5.023' },
    @{ Ref = 'E21'; Value = '  -0.28%  ' },
    @{ Ref = 'D22'; Value = '2.110.87' },
    @{ Ref = 'E22'; Value = '  +0.20%  ' },
    @{ Ref = 'D23'; Value = '10.37' },
    @{ Ref = 'E23'; Value = '  -1.37%  ' },
    @{ Ref = 'D24'; Value = '6.414' },
    @{ Ref = 'E24'; Value = '  -0.46%  ' },
    @{ Ref = 'D25'; Value = '147.59' },
    @{ Ref = 'E25'; Value = '  -1.33%  ' },
    @{ Ref = 'E26'; Value = '  -2.73%  ' },
    @{ Ref = 'D27'; Value = '17.87' },
    @{ Ref = 'E27'; Value = '  -0.39%  ' },
    @{ Ref = 'E28'; Value = '  -0.25%  ' },
    @{ Ref = 'D29'; Value = '113.04' },
    @{ Ref = 'E29'; Value = '  +0.22%  ' },
    @{ Ref = 'D30'; Value = '4.697' },
    @{ Ref = 'E30'; Value = '  +0.66%  ' },
    @{ Ref = 'D31'; Value = '4.667' },
    @{ Ref = 'E31'; Value = '  +0.46%  ' },
    @{ Ref = 'D32'; Value = '0.09225' },
    @{ Ref = 'E32'; Value = '  +2.31%  ' },
    @{ Ref = 'D33'; Value = '0.05109' },
    @{ Ref = 'E33'; Value = '  -0.19%  ' },
    @{ Ref = 'D34'; Value = '3.001' },
    @{ Ref = 'E34'; Value = '  -2.61%  ' },
    @{ Ref = 'D35'; Value = '0.7480' },
    @{ Ref = 'E35'; Value = '  +2.52%  ' },
    @{ Ref = 'D36'; Value = '1.150' },
    @{ Ref = 'E36'; Value = '  +0.28%  ' },
    @{ Ref = 'D37'; Value = '3.272' },
    @{ Ref = 'E37'; Value = '  +7.00%  ' },
    @{ Ref = 'D38'; Value = '2.559' },
    @{ Ref = 'E38'; Value = '  +2.69%  ' },
    @{ Ref = 'D39'; Value = '0.01996' },
    @{ Ref = 'E39'; Value = '  -2.47%  ' },
    @{ Ref = 'D40'; Value = '0.5549' },
    @{ Ref = 'E40'; Value = '  +4.80%  ' },
    @{ Ref = 'D41'; Value = '1.074' },
    @{ Ref = 'E41'; Value = '  +0.04%  ' },
    @{ Ref = 'D42'; Value = '117.61' },
    @{ Ref = 'E42'; Value = '  +2.04%  ' },
    @{ Ref = 'D43'; Value = '6.533' },
    @{ Ref = 'E43'; Value = '  -0.31%  ' },
    @{ Ref = 'D44'; Value = '8.529' },
    @{ Ref = 'E44'; Value = '  +3.76%  ' },
    @{ Ref = 'D45'; Value = '0.1467' },
    @{ Ref = 'E45'; Value = '  +0.27%  ' },
    @{ Ref = 'D46'; Value = '0.4679' },
    @{ Ref = 'E46'; Value = '  +1.89%  ' },
    @{ Ref = 'D47'; Value = '0.9997' },
    @{ Ref = 'E47'; Value = '  -0.06%  ' },
    @{ Ref = 'D48'; Value = '10.03' },
    @{ Ref = 'E48'; Value = '  -0.06%  ' },
    @{ Ref = 'D49'; Value = '1.560' },
    @{ Ref = 'E49'; Value = '  -0.03%  ' },
    @{ Ref = 'D50'; Value = '36.69' },
    @{ Ref = 'E50'; Value = '  +0.73%  ' },
    @{ Ref = 'D51'; Value = '62.97' },
    @{ Ref = 'E51'; Value = '  -1.80%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
